$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Add links for the final assignments (Machine Learning and Many models)
$ws.Range("E5").Value = "[Machine Learning](https://classroom.github.com/a/WRI89Flt)"
$ws.Range("E7").Value = "[Many models](https://classroom.github.com/a/04gGD6TJ)"
